# Dailyweek and Weekly Report 15.02.2025 LeKhanhDuc
# Add a new "I. Status Report" entry in row 16 describing the
# authentication & authorization implementation task.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the formatting of the row above (row 15) onto row 16 so the
# new entry picks up the same borders / fonts / number format / wrap text
# used throughout the "Status Report" table.
$ws.Range("A15:E15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match row 15's taller row height for the wrapped notes text.
$ws.Rows.Item(16).RowHeight = 60

# Fill in the new task details.
$ws.Range("B16").Value = "Implement authentication & authorization for backend security."
$ws.Range("C16").Value = "Le Khanh Duc"
$ws.Range("D16").Value = "'09/02/2025"
$ws.Range("E16").Value = "Implement authentication & authorization in ASP.NET Core. Focus on security best practices, token-based authentication (JWT), and role-based access control (RBAC) to ensure secure API access"

# Move the active selection to the new row, as in the edited workbook.
$ws.Range("E16").Select()
